# Update the FAST_holdings model holdings workbook:
#  - bump the "as of" date in the confidential disclosure note from
#    2021-04-21 to 2021-04-22
#  - refresh the Weight (D) and Percent Change (E) figures for rows 2-10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect (using the workbook's known password)
# before writing, then restore protection afterwards.
$ws.Unprotect("D382")

$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-22 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.1014998070682455
$ws.Range("E2").Value = 0.003915115751248255

$ws.Range("D3").Value = 0.1090582860807263
$ws.Range("E3").Value = -0.006729709346173895

$ws.Range("D4").Value = 0.1178085099972805
$ws.Range("E4").Value = -0.007378335949764669

$ws.Range("D5").Value = 0.1377635045767312
$ws.Range("E5").Value = -0.004964766175528545

$ws.Range("D6").Value = 0.1352596397728285
$ws.Range("E6").Value = -0.005517241379310311

$ws.Range("D7").Value = 0.1424467732795342
$ws.Range("E7").Value = -0.009375298957237321

$ws.Range("D8").Value = 0.1277224129496254
$ws.Range("E8").Value = -0.004966403739409819

$ws.Range("D9").Value = 0.1284410662750285
$ws.Range("E9").Value = -0.0008254056236028928

$ws.Range("E10").Value = -0.00471181964013867

$ws.Protect("D382")
